$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChartSettings")

# --- Row 2: rename the ID and add the missing Date_Change, fix U2 ---
$ws.Range("A2").Value = "ExpectationRanges.01"
$ws.Range("B2").Value = "9/6/2024"
$ws.Range("U2").Value = 0.12

# --- Row 3: new supplemented row ---
$ws.Range("A3").Value = "ExpectationRanges.02"
$ws.Range("B3").Value = "7/5/2024"
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("G3").Value = $ws.Range("G2").Value2
$ws.Range("H3").Value = $ws.Range("H2").Value2
$ws.Range("O3").Value = $ws.Range("O2").Value2
$ws.Range("P3").Value = $ws.Range("P2").Value2
$ws.Range("R3").Value = 20
$ws.Range("S3").Value = 20
$ws.Range("T3").Value = 0.12
$ws.Range("U3").Value = 0.2
$ws.Range("V3").Value = 0.3
$ws.Range("W3").Value = 0.1
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 400
$ws.Range("Z3").Value = 20
$ws.Range("AA3").Value = 5
$ws.Range("AB3").Value = 50

# --- Row 4: new supplemented row ---
$ws.Range("A4").Value = "ExpectationRanges.03"
$ws.Range("B4").Value = "9/6/2024"
$ws.Range("D4").Value = $ws.Range("D2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2
$ws.Range("G4").Value = $ws.Range("G2").Value2
$ws.Range("H4").Value = $ws.Range("H2").Value2
$ws.Range("O4").Value = $ws.Range("O2").Value2
$ws.Range("P4").Value = $ws.Range("P2").Value2
$ws.Range("R4").Value = 20
$ws.Range("S4").Value = 20
$ws.Range("T4").Value = 0.12
$ws.Range("U4").Value = 0.25
$ws.Range("V4").Value = 0.3
$ws.Range("W4").Value = 0.1
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 400
$ws.Range("Z4").Value = 20
$ws.Range("AA4").Value = 5
$ws.Range("AB4").Value = 50

# --- Row 5: new supplemented row ---
$ws.Range("A5").Value = "ExpectationRanges.04"
$ws.Range("B5").Value = "9/6/2024"
$ws.Range("D5").Value = $ws.Range("D2").Value2
$ws.Range("E5").Value = $ws.Range("E2").Value2
$ws.Range("G5").Value = $ws.Range("G2").Value2
$ws.Range("H5").Value = $ws.Range("H2").Value2
$ws.Range("O5").Value = $ws.Range("O2").Value2
$ws.Range("P5").Value = $ws.Range("P2").Value2
$ws.Range("R5").Value = 20
$ws.Range("S5").Value = 20
$ws.Range("T5").Value = 0.12
$ws.Range("U5").Value = 0.3
$ws.Range("V5").Value = 0.3
$ws.Range("W5").Value = 0.1
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 400
$ws.Range("Z5").Value = 20
$ws.Range("AA5").Value = 5
$ws.Range("AB5").Value = 50

$ws.Range("U5").Select()
